$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (dates as Excel serial numbers, error counts)
$dates  = @(45957, 45959, 45966, 45967, 45964, 45968, 45958, 45960, 45965, 45961)
$counts = @(64, 80, 124, 58, 110, 63, 82, 86, 112, 108)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $counts[$i]
}

# Copy formatting from an already-formatted date cell (A2) down into the new rows
# so the new cells pick up the existing date-number-format style instead of
# creating a brand new style entry.
$ws.Range("A2").Copy()
$ws.Range("A7:A11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update selection to match the new range
$ws.Range("A2:B11").Select()
